# Remove the two example/reference tables from the template body.
# (The short "ada/Ada" spell-check demo table and the longer
#  1.1.x requirements-checklist demo table.) The surrounding
# paragraphs are left untouched.
$d = $word.ActiveDocument

# Delete from the end backwards so earlier indices stay valid.
for ($i = $d.Tables.Count; $i -ge 1; $i--) {
    $d.Tables.Item($i).Delete()
}
